$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 12.90258951017826
$ws.Cells.Item(2, 3).Value = 8.765811693360169
$ws.Cells.Item(2, 5).Value = 15.3704130336887
$ws.Cells.Item(2, 6).Value = 39.46161740653699
$ws.Cells.Item(2, 7).Value = 3.6727954155914
$ws.Cells.Item(2, 9).Value = 26.5361954024308
$ws.Cells.Item(2, 10).Value = 8.350327027231986
$ws.Cells.Item(2, 11).Value = 9.56381890483544
$ws.Cells.Item(2, 12).Value = 12.45658811191793
$ws.Cells.Item(2, 15).Value = 26.33011577458971
$ws.Cells.Item(3, 2).Value = 12.63982079858845
$ws.Cells.Item(3, 3).Value = 8.76215732473
$ws.Cells.Item(3, 5).Value = 15.36268560497533
$ws.Cells.Item(3, 6).Value = 39.51074005399009
$ws.Cells.Item(3, 7).Value = 3.67450500554647
$ws.Cells.Item(3, 9).Value = 26.64713912378761
$ws.Cells.Item(3, 10).Value = 8.346808597540042
$ws.Cells.Item(3, 11).Value = 9.382713505093147
$ws.Cells.Item(3, 12).Value = 12.43058558535039
$ws.Cells.Item(3, 15).Value = 26.43619226736916
$ws.Cells.Item(4, 2).Value = 12.4773962439341
$ws.Cells.Item(4, 3).Value = 8.76011552084781
$ws.Cells.Item(4, 5).Value = 15.36019872198642
$ws.Cells.Item(4, 6).Value = 39.54968857524139
$ws.Cells.Item(4, 7).Value = 3.675610278513336
$ws.Cells.Item(4, 9).Value = 26.71983979859319
$ws.Cells.Item(4, 10).Value = 8.344741556295403
$ws.Cells.Item(4, 11).Value = 9.270802397242415
$ws.Cells.Item(4, 12).Value = 12.4163171672956
$ws.Cells.Item(4, 15).Value = 26.50628634793635
$ws.Cells.Item(5, 2).Value = 12.41102263467561
$ws.Cells.Item(5, 3).Value = 8.759334664144539
$ws.Cells.Item(5, 5).Value = 15.35975524528305
$ws.Cells.Item(5, 6).Value = 39.56776805884027
$ws.Cells.Item(5, 7).Value = 3.676074704634876
$ws.Cells.Item(5, 9).Value = 26.7506180181749
$ws.Cells.Item(5, 10).Value = 8.343922888099225
$ws.Cells.Item(5, 11).Value = 9.225075879816231
$ws.Cells.Item(5, 12).Value = 12.41093350547315
$ws.Cells.Item(5, 15).Value = 26.53609729564045
$ws.Cells.Item(6, 2).Value = 12.39999284297545
$ws.Cells.Item(6, 3).Value = 8.759208107267304
$ws.Cells.Item(6, 5).Value = 15.35971608957903
$ws.Cells.Item(6, 6).Value = 39.57090340838003
$ws.Cells.Item(6, 7).Value = 3.676152670178733
$ws.Cells.Item(6, 9).Value = 26.75579830981956
$ws.Cells.Item(6, 10).Value = 8.343788384446594
$ws.Cells.Item(6, 11).Value = 9.21747736460501
$ws.Cells.Item(6, 12).Value = 12.4100656864002
$ws.Cells.Item(6, 15).Value = 26.54112267404948
$ws.Cells.Item(7, 2).Value = 12.47650173403518
$ws.Cells.Item(7, 3).Value = 8.760104782121914
$ws.Cells.Item(7, 5).Value = 15.36019043066204
$ws.Cells.Item(7, 6).Value = 39.54992346592975
$ws.Cells.Item(7, 7).Value = 3.675616485109508
$ws.Cells.Item(7, 9).Value = 26.72025021927842
$ws.Cells.Item(7, 10).Value = 8.34473041934605
$ws.Cells.Item(7, 11).Value = 9.270186130234379
$ws.Cells.Item(7, 12).Value = 12.41624281178436
$ws.Cells.Item(7, 15).Value = 26.50668334072545
$ws.Cells.Item(8, 2).Value = 12.81226381200265
$ws.Cells.Item(8, 3).Value = 8.764509983572252
$ws.Cells.Item(8, 5).Value = 15.36728140535745
$ws.Cells.Item(8, 6).Value = 39.47673021596179
$ws.Cells.Item(8, 7).Value = 3.673373372913299
$ws.Cells.Item(8, 9).Value = 26.57349819153129
$ws.Cells.Item(8, 10).Value = 8.349094635771017
$ws.Cells.Item(8, 11).Value = 9.501555275493599
$ws.Cells.Item(8, 12).Value = 12.44727269320262
$ws.Cells.Item(8, 15).Value = 26.36566042281178
$ws.Cells.Item(9, 2).Value = 13.45836709505537
$ws.Cells.Item(9, 3).Value = 8.774734248425704
$ws.Cells.Item(9, 5).Value = 15.39899580354193
$ws.Cells.Item(9, 6).Value = 39.40297621254455
$ws.Cells.Item(9, 7).Value = 3.669413660258605
$ws.Cells.Item(9, 9).Value = 26.32205378174618
$ws.Cells.Item(9, 10).Value = 8.358384274812641
$ws.Cells.Item(9, 11).Value = 9.947195627131135
$ws.Cells.Item(9, 12).Value = 12.52139979916696
$ws.Cells.Item(9, 15).Value = 26.12853024589825
$ws.Cells.Item(10, 2).Value = 13.92079574732447
$ws.Cells.Item(10, 3).Value = 8.783190692005906
$ws.Cells.Item(10, 5).Value = 15.43299920550884
$ws.Cells.Item(10, 6).Value = 39.39136150489407
$ws.Cells.Item(10, 7).Value = 3.666769349698653
$ws.Cells.Item(10, 9).Value = 26.15945403648488
$ws.Cells.Item(10, 10).Value = 8.365644102494374
$ws.Cells.Item(10, 11).Value = 10.26657523131038
$ws.Cells.Item(10, 12).Value = 12.58368560823891
$ws.Cells.Item(10, 15).Value = 25.97839214309356
$ws.Cells.Item(11, 2).Value = 14.12757740055507
$ws.Cells.Item(11, 3).Value = 8.787238199909526
$ws.Cells.Item(11, 5).Value = 15.45075437707397
$ws.Cells.Item(11, 6).Value = 39.39531179593273
$ws.Cells.Item(11, 7).Value = 3.665623317827438
$ws.Cells.Item(11, 9).Value = 26.09028761069463
$ws.Cells.Item(11, 10).Value = 8.36903879806181
$ws.Cells.Item(11, 11).Value = 10.40951308248613
$ws.Cells.Item(11, 12).Value = 12.61365684661879
$ws.Cells.Item(11, 15).Value = 25.9153313110757
$ws.Cells.Item(12, 2).Value = 14.20529420850412
$ws.Cells.Item(12, 3).Value = 8.788799317163992
$ws.Cells.Item(12, 5).Value = 15.45780298257003
$ws.Cells.Item(12, 6).Value = 39.39813302147001
$ws.Cells.Item(12, 7).Value = 3.665197480232642
$ws.Cells.Item(12, 9).Value = 26.0647865020483
$ws.Cells.Item(12, 10).Value = 8.370337311871893
$ws.Cells.Item(12, 11).Value = 10.4632545554808
$ws.Cells.Item(12, 12).Value = 12.62523577199169
$ws.Cells.Item(12, 15).Value = 25.89220611855012
$ws.Cells.Item(13, 2).Value = 14.18858370668239
$ws.Cells.Item(13, 3).Value = 8.78846184589869
$ws.Cells.Item(13, 5).Value = 15.45627054224488
$ws.Cells.Item(13, 6).Value = 39.39746652482999
$ws.Cells.Item(13, 7).Value = 3.665288830578938
$ws.Cells.Item(13, 9).Value = 26.07024789649609
$ws.Cells.Item(13, 10).Value = 8.370057078587067
$ws.Cells.Item(13, 11).Value = 10.45169827086404
$ws.Cells.Item(13, 12).Value = 12.62273193500802
$ws.Cells.Item(13, 15).Value = 25.89715296016968
$ws.Cells.Item(14, 2).Value = 14.13398336878862
$ws.Cells.Item(14, 3).Value = 8.787366067019956
$ws.Cells.Item(14, 5).Value = 15.45132777676773
$ws.Cells.Item(14, 6).Value = 39.3955173556752
$ws.Cells.Item(14, 7).Value = 3.665588121002057
$ws.Cells.Item(14, 9).Value = 26.08817577409135
$ws.Cells.Item(14, 10).Value = 8.369145368075358
$ws.Cells.Item(14, 11).Value = 10.41394242721775
$ws.Cells.Item(14, 12).Value = 12.61460488787368
$ws.Cells.Item(14, 15).Value = 25.91341365669565
$ws.Cells.Item(15, 2).Value = 14.10046058812011
$ws.Cells.Item(15, 3).Value = 8.786698558136116
$ws.Cells.Item(15, 5).Value = 15.44834241348843
$ws.Cells.Item(15, 6).Value = 39.39449594160066
$ws.Cells.Item(15, 7).Value = 3.66577250391156
$ws.Cells.Item(15, 9).Value = 26.09924707904856
$ws.Cells.Item(15, 10).Value = 8.368588607089121
$ws.Cells.Item(15, 11).Value = 10.39076423552381
$ws.Cells.Item(15, 12).Value = 12.60965653944958
$ws.Cells.Item(15, 15).Value = 25.92347211257658
$ws.Cells.Item(16, 2).Value = 13.90720377493565
$ws.Cells.Item(16, 3).Value = 8.782930177784575
$ws.Cells.Item(16, 5).Value = 15.4318845667513
$ws.Cells.Item(16, 6).Value = 39.3912889621068
$ws.Cells.Item(16, 7).Value = 3.666845386751537
$ws.Cells.Item(16, 9).Value = 26.16407086146099
$ws.Cells.Item(16, 10).Value = 8.365424085550883
$ws.Cells.Item(16, 11).Value = 10.25718244943611
$ws.Cells.Item(16, 12).Value = 12.58175935776909
$ws.Cells.Item(16, 15).Value = 25.98261883241341
$ws.Cells.Item(17, 2).Value = 13.78767736084741
$ws.Cells.Item(17, 3).Value = 8.780669474921366
$ws.Cells.Item(17, 5).Value = 15.42237134876687
$ws.Cells.Item(17, 6).Value = 39.39168512806954
$ws.Cells.Item(17, 7).Value = 3.667518105500418
$ws.Cells.Item(17, 9).Value = 26.20506803236177
$ws.Cells.Item(17, 10).Value = 8.363506191869613
$ws.Cells.Item(17, 11).Value = 10.1745974438984
$ws.Cells.Item(17, 12).Value = 12.56506049567223
$ws.Cells.Item(17, 15).Value = 26.02024592817591
$ws.Cells.Item(18, 2).Value = 13.71859604003634
$ws.Cells.Item(18, 3).Value = 8.77938808890422
$ws.Cells.Item(18, 5).Value = 15.4171150867921
$ws.Cells.Item(18, 6).Value = 39.39278233582862
$ws.Cells.Item(18, 7).Value = 3.667910391507417
$ws.Cells.Item(18, 9).Value = 26.2291003920303
$ws.Cells.Item(18, 10).Value = 8.362411759353019
$ws.Cells.Item(18, 11).Value = 10.1268782866441
$ws.Cells.Item(18, 12).Value = 12.55561021726667
$ws.Cells.Item(18, 15).Value = 26.04238095515007
$ws.Cells.Item(19, 2).Value = 13.69515139339704
$ws.Cells.Item(19, 3).Value = 8.778957494753948
$ws.Cells.Item(19, 5).Value = 15.41537252710046
$ws.Cells.Item(19, 6).Value = 39.39330320268728
$ws.Cells.Item(19, 7).Value = 3.668044133903807
$ws.Cells.Item(19, 9).Value = 26.23731493591138
$ws.Cells.Item(19, 10).Value = 8.362042702854238
$ws.Cells.Item(19, 11).Value = 10.11068539411916
$ws.Cells.Item(19, 12).Value = 12.55243722188811
$ws.Cells.Item(19, 15).Value = 26.04996011167111
$ws.Cells.Item(20, 2).Value = 13.80043614208367
$ws.Cells.Item(20, 3).Value = 8.780908176591868
$ws.Cells.Item(20, 5).Value = 15.42336177100836
$ws.Cells.Item(20, 6).Value = 39.39155299677725
$ws.Cells.Item(20, 7).Value = 3.667445939403188
$ws.Cells.Item(20, 9).Value = 26.20065704423864
$ws.Cells.Item(20, 10).Value = 8.363709457330005
$ws.Cells.Item(20, 11).Value = 10.18341172468299
$ws.Cells.Item(20, 12).Value = 12.56682217371821
$ws.Cells.Item(20, 15).Value = 26.01618943450392
$ws.Cells.Item(21, 2).Value = 14.15003729259814
$ws.Cells.Item(21, 3).Value = 8.787687156186836
$ws.Cells.Item(21, 5).Value = 15.45277079453204
$ws.Cells.Item(21, 6).Value = 39.39605392906342
$ws.Cells.Item(21, 7).Value = 3.665499991548653
$ws.Cells.Item(21, 9).Value = 26.08289117314727
$ws.Cells.Item(21, 10).Value = 8.369412808195834
$ws.Cells.Item(21, 11).Value = 10.42504307886327
$ws.Cells.Item(21, 12).Value = 12.61698581959219
$ws.Cells.Item(21, 15).Value = 25.90861700578808
$ws.Cells.Item(22, 2).Value = 14.37507066209552
$ws.Cells.Item(22, 3).Value = 8.792283024353589
$ws.Cells.Item(22, 5).Value = 15.47388435463686
$ws.Cells.Item(22, 6).Value = 39.40671881237529
$ws.Cells.Item(22, 7).Value = 3.664275629164995
$ws.Cells.Item(22, 9).Value = 26.00995057467105
$ws.Cells.Item(22, 10).Value = 8.373216034603429
$ws.Cells.Item(22, 11).Value = 10.58069298475222
$ws.Cells.Item(22, 12).Value = 12.65110523981329
$ws.Cells.Item(22, 15).Value = 25.84271110981158
$ws.Cells.Item(23, 2).Value = 14.25530454846908
$ws.Cells.Item(23, 3).Value = 8.789815132646853
$ws.Cells.Item(23, 5).Value = 15.46244373503498
$ws.Cells.Item(23, 6).Value = 39.40032114511806
$ws.Cells.Item(23, 7).Value = 3.664924767709761
$ws.Cells.Item(23, 9).Value = 26.04851182457988
$ws.Cells.Item(23, 10).Value = 8.371179328295375
$ws.Cells.Item(23, 11).Value = 10.49784258508166
$ws.Cells.Item(23, 12).Value = 12.63277497933771
$ws.Cells.Item(23, 15).Value = 25.87748336175551
$ws.Cells.Item(24, 2).Value = 13.79466902530559
$ws.Cells.Item(24, 3).Value = 8.78080020246837
$ws.Cells.Item(24, 5).Value = 15.42291333732613
$ws.Cells.Item(24, 6).Value = 39.39161002489676
$ws.Cells.Item(24, 7).Value = 3.667478548464035
$ws.Cells.Item(24, 9).Value = 26.20264981119392
$ws.Cells.Item(24, 10).Value = 8.363617535482891
$ws.Cells.Item(24, 11).Value = 10.17942753200115
$ws.Cells.Item(24, 12).Value = 12.5660252510235
$ws.Cells.Item(24, 15).Value = 26.01802180947437
$ws.Cells.Item(25, 2).Value = 13.28540388882403
$ws.Cells.Item(25, 3).Value = 8.771800816226436
$ws.Cells.Item(25, 5).Value = 15.38852452798935
$ws.Cells.Item(25, 6).Value = 39.415451534796
$ws.Cells.Item(25, 7).Value = 3.670438152451468
$ws.Cells.Item(25, 9).Value = 26.38618785131153
$ws.Cells.Item(25, 10).Value = 8.355794019533679
$ws.Cells.Item(25, 11).Value = 9.827832352540645
$ws.Cells.Item(25, 12).Value = 12.49995168929223
$ws.Cells.Item(25, 15).Value = 26.18845518802883
